$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these price cells to remain Text (matches source inlineStr values)
# since their new values look like plain numbers and would otherwise be
# auto-converted to numeric by Excel's input parsing.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.437.73'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '3.139.31'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '608.93'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = '143.73'
$ws.Range("E6").Value = '  -2.44%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.136.25'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").Value = '5.38'
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("D12").Value = '0.471'
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").Value = '35.38'
$ws.Range("E14").Value = '  -0.74%  '
$ws.Range("D15").Value = '3.651.09'
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("E16").Value = '  +2.54%  '
$ws.Range("D17").Value = '64.391.92'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").Value = '3.126.58'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = '6.86'
$ws.Range("E19").Value = '  -1.07%  '
$ws.Range("D20").Value = '476.72'
$ws.Range("E20").Value = '  -0.83%  '
$ws.Range("D21").Value = '14.85'
$ws.Range("E21").Value = '  +0.62%  '
$ws.Range("D22").Value = '0.718'
$ws.Range("E22").Value = '  +1.02%  '
$ws.Range("D23").Value = '7.79'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '85.77'
$ws.Range("E24").Value = '  +2.44%  '
$ws.Range("D25").Value = '13.52'
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("D28").Value = '8.46'
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '7.34'
$ws.Range("E29").Value = '  +7.41%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '0.115'
$ws.Range("E30").Value = '  +1.83%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '2.06'
$ws.Range("E31").Value = '  -6.22%  '
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").Value = '26.70'
$ws.Range("E33").Value = '  +1.78%  '
$ws.Range("E34").Value = '  -3.77%  '
$ws.Range("E35").Value = '  +0.29%  '
$ws.Range("D36").Value = '5.97'
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").Value = '52.69'
$ws.Range("D38").Value = '0.0₃0739'
$ws.Range("E38").Value = '  +0.96%  '
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Value = '2.99'
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '448.36'
$ws.Range("E40").Value = '  -1.83%  '
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = '8.32'
$ws.Range("E43").Value = '  -1.42%  '
$ws.Range("D44").Value = '2.883.31'
$ws.Range("E44").Value = '  +1.10%  '
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("E46").Value = '  -1.65%  '
$ws.Range("E47").Value = '  +4.55%  '
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '120.71'
$ws.Range("E51").Value = '  +0.79%  '
